# Auto-generated Excel COM-interop script
# Applies numeric cell-value corrections to the Sagittarius_Profits workbook
# as captured by the authoritative OOXML diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 188.33333
$ws.Range("I5").Value = 33.25
$ws.Range("J5").Value = 312.4
$ws.Range("K5").Value = 33.25
$ws.Range("L5").Value = 312.4
$ws.Range("M5").Value = 81.75
$ws.Range("N5").Value = -542.4
$ws.Range("H8").Value = 1628.5714
$ws.Range("I8").Value = 20.333334
$ws.Range("J8").Value = 2834.75
$ws.Range("K8").Value = 61.000002
$ws.Range("L8").Value = 8504.25
$ws.Range("M8").Value = 77.99999800000001
$ws.Range("N8").Value = -8782.25
$ws.Range("H12").Value = 3197.5
$ws.Range("I12").Value = 396
$ws.Range("K12").Value = 396
$ws.Range("M12").Value = -226
$ws.Range("H19").Value = 732.3333
$ws.Range("I19").Value = 465.66666
$ws.Range("J19").Value = 865.6667
$ws.Range("K19").Value = 465.66666
$ws.Range("L19").Value = 865.6667
$ws.Range("M19").Value = -290.66666
$ws.Range("N19").Value = -1215.6667
$ws.Range("H28").Value = 494
$ws.Range("I28").Value = 447.14285
$ws.Range("K28").Value = 447.14285
$ws.Range("M28").Value = 37.85714999999999
$ws.Range("H98").Value = 1324
$ws.Range("I98").Value = 1324
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1324
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 174
$ws.Range("N98").Value = ""
$ws.Range("H122").Value = 1324
$ws.Range("I122").Value = 1324
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3972
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1522
$ws.Range("N122").Value = ""
$ws.Range("H125").Value = 23666.666
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 23666.666
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 212999.994
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -217919.994
$ws.Range("H127").Value = 1792.2222
$ws.Range("I127").Value = 1819.1428
$ws.Range("J127").Value = 1698
$ws.Range("K127").Value = 5457.428400000001
$ws.Range("L127").Value = 5094
$ws.Range("M127").Value = -497.4284000000007
$ws.Range("N127").Value = -15014
$ws.Range("H132").Value = 2997.9167
$ws.Range("I132").Value = 2997.7273
$ws.Range("K132").Value = 8993.1819
$ws.Range("M132").Value = -6463.1819

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 4826.6665
$ws.Range("J46").Value = 4826.6665
$ws.Range("L46").Value = 4826.6665
$ws.Range("N46").Value = -5464.6665
$ws.Range("H61").Value = 3747.5
$ws.Range("I61").Value = 3330
$ws.Range("K61").Value = 3330
$ws.Range("M61").Value = -3118
$ws.Range("H122").Value = 1133.3334
$ws.Range("J122").Value = 1150
$ws.Range("L122").Value = 3450
$ws.Range("N122").Value = -8350
$ws.Range("H136").Value = 3747.5
$ws.Range("I136").Value = 3330
$ws.Range("K136").Value = 9990
$ws.Range("M136").Value = -7440

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 10000
$ws.Range("I96").Value = 10000
$ws.Range("K96").Value = 10000
$ws.Range("M96").Value = -7254
$ws.Range("H107").Value = 1269.6428
$ws.Range("I107").Value = 1142.5555
$ws.Range("K107").Value = 1142.5555
$ws.Range("M107").Value = 777.4445000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 370.57144
$ws.Range("I22").Value = 239.8
$ws.Range("K22").Value = 239.8
$ws.Range("M22").Value = 110.2
$ws.Range("H58").Value = 2349.5715
$ws.Range("I58").Value = 1948.5
$ws.Range("J58").Value = 2884.3333
$ws.Range("K58").Value = 1948.5
$ws.Range("L58").Value = 2884.3333
$ws.Range("M58").Value = -1745.5
$ws.Range("N58").Value = -3290.3333
$ws.Range("H99").Value = 1258
$ws.Range("I99").Value = 980.4
$ws.Range("J99").Value = 2183.3333
$ws.Range("K99").Value = 980.4
$ws.Range("L99").Value = 2183.3333
$ws.Range("M99").Value = 517.6
$ws.Range("N99").Value = -5179.3333
$ws.Range("H122").Value = 1305.25
$ws.Range("I122").Value = 1177.4286
$ws.Range("K122").Value = 3532.2858
$ws.Range("M122").Value = -1082.2858
$ws.Range("H126").Value = 1258
$ws.Range("I126").Value = 980.4
$ws.Range("J126").Value = 2183.3333
$ws.Range("K126").Value = 2941.2
$ws.Range("L126").Value = 6549.999899999999
$ws.Range("M126").Value = -471.1999999999998
$ws.Range("N126").Value = -11489.9999
$ws.Range("H136").Value = 2349.5715
$ws.Range("I136").Value = 1948.5
$ws.Range("J136").Value = 2884.3333
$ws.Range("K136").Value = 5845.5
$ws.Range("L136").Value = 8652.999899999999
$ws.Range("M136").Value = -3295.5
$ws.Range("N136").Value = -13752.9999
$ws.Range("H141").Value = 38037.145
$ws.Range("J141").Value = 38037.145
$ws.Range("L141").Value = 38037.145
$ws.Range("N141").Value = -48397.145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 999
$ws.Range("J2").Value = 999
$ws.Range("L2").Value = 5994
$ws.Range("N2").Value = -6220
$ws.Range("H119").Value = 3500
$ws.Range("I119").Value = 3500
$ws.Range("K119").Value = 10500
$ws.Range("M119").Value = -5662
$ws.Range("H131").Value = 502069.9
$ws.Range("J131").Value = 911527.0600000001
$ws.Range("L131").Value = 2734581.18
$ws.Range("N131").Value = -2744661.18

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4937.75
$ws.Range("I70").Value = 4917
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4917
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4647
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 4937.75
$ws.Range("I73").Value = 4917
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4917
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -3981
$ws.Range("N73").Value = -6872
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = ""
$ws.Range("H102").Value = 1885.2727
$ws.Range("I102").Value = 1687.75
$ws.Range("J102").Value = 1998.1428
$ws.Range("K102").Value = 1687.75
$ws.Range("L102").Value = 1998.1428
$ws.Range("M102").Value = -65.75
$ws.Range("N102").Value = -5242.1428
$ws.Range("H122").Value = 2582.9
$ws.Range("I122").Value = 2715.1765
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 8145.529500000001
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -5695.529500000001
$ws.Range("N122").Value = -10400.0002
$ws.Range("H132").Value = 7987.727
$ws.Range("I132").Value = 8101.7646
$ws.Range("K132").Value = 24305.2938
$ws.Range("M132").Value = -21775.2938
$ws.Range("H134").Value = 36221
$ws.Range("J134").Value = 36221
$ws.Range("L134").Value = 108663
$ws.Range("N134").Value = -113733

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9786.066000000001
$ws.Range("I7").Value = 11250
$ws.Range("K7").Value = 11250
$ws.Range("M7").Value = -11138
$ws.Range("H40").Value = 4499.5454
$ws.Range("I40").Value = 3936.875
$ws.Range("K40").Value = 3936.875
$ws.Range("M40").Value = -3800.875
$ws.Range("H46").Value = 1527.9231
$ws.Range("I46").Value = 1218.1
$ws.Range("J46").Value = 1854.0526
$ws.Range("K46").Value = 1218.1
$ws.Range("L46").Value = 1854.0526
$ws.Range("M46").Value = -1030.1
$ws.Range("N46").Value = -2230.0526
$ws.Range("H122").Value = 8471.130999999999
$ws.Range("I122").Value = 9045.1
$ws.Range("K122").Value = 27135.3
$ws.Range("M122").Value = -24685.3
$ws.Range("H126").Value = 9786.066000000001
$ws.Range("I126").Value = 11250
$ws.Range("K126").Value = 33750
$ws.Range("M126").Value = -31280
$ws.Range("H132").Value = 3766.3333
$ws.Range("I132").Value = 3766.3333
$ws.Range("K132").Value = 11298.9999
$ws.Range("M132").Value = -8768.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("M29").Value = -710
$ws.Range("H74").Value = 19797.25
$ws.Range("J74").Value = 20303.834
$ws.Range("L74").Value = 20303.834
$ws.Range("N74").Value = -22175.834
$ws.Range("H77").Value = 19797.25
$ws.Range("J77").Value = 20303.834
$ws.Range("L77").Value = 60911.50199999999
$ws.Range("N77").Value = -70271.50199999999
$ws.Range("H81").Value = 1002069.7
$ws.Range("I81").Value = 1549.5
$ws.Range("J81").Value = 1669083.1
$ws.Range("K81").Value = 3099
$ws.Range("L81").Value = 3338166.2
$ws.Range("M81").Value = -2038
$ws.Range("N81").Value = -3340288.2
$ws.Range("H84").Value = 1002069.7
$ws.Range("I84").Value = 1549.5
$ws.Range("J84").Value = 1669083.1
$ws.Range("K84").Value = 15495
$ws.Range("L84").Value = 16690831
$ws.Range("M84").Value = -10191
$ws.Range("N84").Value = -16701439
$ws.Range("H122").Value = 938.6667
$ws.Range("I122").Value = 556
$ws.Range("K122").Value = 1668
$ws.Range("M122").Value = 782
$ws.Range("H126").Value = 2638.2222
$ws.Range("I126").Value = 2192
$ws.Range("K126").Value = 6576
$ws.Range("M126").Value = -4106
